$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add missing additional-start-date (column C) values for several policy rows,
# matching the date format already used by the surrounding C/D columns.
$ws.Range("C8").Value = 43938
$ws.Range("C8").NumberFormat = "YYYY-MM-DD"

$ws.Range("C9").Value = 43907
$ws.Range("C9").NumberFormat = "YYYY-MM-DD"

$ws.Range("C15").Value = 43907
$ws.Range("C15").NumberFormat = "YYYY-MM-DD"

$ws.Range("C20").Value = 43915
$ws.Range("C20").NumberFormat = "YYYY-MM-DD"

$ws.Range("C24").Value = 43915
$ws.Range("C24").NumberFormat = "YYYY-MM-DD"

# Reflect the author's final view/selection state in the saved workbook.
$ws.Range("C24").Select()
